$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 cell text to the new data path (selectDeviceOptionsAndProceedToCheckout instead of placeOrderAndProceedToCheckout)
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/selectDeviceOptionsAndProceedToCheckout-test-data"

# Delete columns E and F entirely (drops link_byPrice_internalRoleLinkName / link_byPrice_nthChild columns)
$ws.Range("E1:F2").EntireColumn.Delete()

# Widen column A from 73 to 82 (ColumnWidth applies Excel's +5/6 character padding, so
# subtract that offset to land exactly on a stored width of 82)
$ws.Columns.Item(1).ColumnWidth = 81.16666666666667
